$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-15 03:57:22"
$wsZh.Range("H2").Value = "2016-03-15 03:58:04"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-15 03:57:30"
$wsDe.Range("H2").Value = "2016-03-15 03:58:17"
